$d = $word.ActiveDocument

# --- 1. Insert the new "File CourseTest.java" glossary entry -------------
# (itIsResultThatStudentStudentIsAbleToTakeTheCourse) right after the
# "thatStudentWantsToTakeTheCourse(): ..." bullet (paragraph 11, 1-indexed,
# in the original document).
$stPara = $d.Paragraphs(11)
$null = $stPara.Range.InsertParagraphAfter()
$d.Paragraphs(12).Range.Text = "itIsResultThatStudentStudentIsAbleToTakeTheCourse(): a method to check if the Student is put in the class correctly."

# --- 2. Insert the two new "GetStudents.java" glossary entries -----------
# (checkCourse / getApproval) right after the
# "GetStudents – A class with methods  that handles different student
# cases" bullet. One extra paragraph was inserted above, so this bullet is
# now paragraph 18 (was 17 before step 1).
$gsPara = $d.Paragraphs(18)
$null = $gsPara.Range.InsertParagraphAfter()
$d.Paragraphs(19).Range.Text = "checkCourse(): a method to split input from CSV, add data to the correct vector. "
$null = $d.Paragraphs(19).Range.InsertParagraphAfter()
$d.Paragraphs(20).Range.Text = "getApproval(): a method to check if the student meets the course requirement."

# --- 3. Drop the redundant " in the GetStudent class" phrase -------------
# from the remaining GetStudents.java entries (getHours, getNumberOfCourses,
# compare, getName, getResult, addTakenCourses).
$null = $d.Content.Find.Execute(" in the GetStudent class", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
